$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.572.98'
$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").Value = '2.461.81'
$ws.Range("E3").Value = '  +0.93%  '

$ws.Range("E4").Value = '  -1.10%  '

$ws.Range("D5").Value = "'314.81"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.76%  '

$ws.Range("D6").Value = "'91.65"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.45%  '

$ws.Range("D7").Value = "'0.548"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.52%  '

$ws.Range("E8").Value = '  -1.06%  '

$ws.Range("D9").Value = "'0.511"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +6.35%  '

$ws.Range("D10").Value = "'32.55"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.11%  '

$ws.Range("E11").Value = '  +4.22%  '

$ws.Range("E12").Value = '  +1.57%  '

$ws.Range("D13").Value = '2.839.37'
$ws.Range("E13").Value = '  +0.96%  '

$ws.Range("D14").Value = "'6.85"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.86%  '

$ws.Range("D15").Value = "'15.83"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +6.03%  '

$ws.Range("D16").Value = '2.445.41'
$ws.Range("E16").Value = '  +0.95%  '

$ws.Range("D17").Value = "'0.772"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +3.20%  '

$ws.Range("D18").Value = '41.571.60'
$ws.Range("E18").Value = '  +1.12%  '

$ws.Range("D19").Value = "'6.48"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +5.78%  '

$ws.Range("D20").Value = '0.0₃0937'
$ws.Range("E20").Value = '  +4.34%  '

$ws.Range("D21").Value = "'70.93"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.16%  '

$ws.Range("D22").Value = "'11.36"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +6.68%  '

$ws.Range("D23").Value = "'237.11"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.07%  '

$ws.Range("E24").Value = '  +2.26%  '

$ws.Range("B25").Value = 'ImmutableX'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D25").Value = "'1.91"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.90%  '

$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.24%  '

$ws.Range("D27").Value = "'24.31"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.71%  '

$ws.Range("D28").Value = "'2.26"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.58%  '

$ws.Range("D29").Value = "'9.66"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.59%  '

$ws.Range("D30").Value = "'35.19"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.32%  '

$ws.Range("D31").Value = "'156.04"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.30%  '

$ws.Range("D32").Value = "'5.45"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +4.52%  '

$ws.Range("D33").Value = "'2.58"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.25%  '

$ws.Range("D34").Value = "'0.0759"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.50%  '

$ws.Range("D35").Value = "'17.45"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.78%  '

$ws.Range("D36").Value = "'2.41"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.29%  '

$ws.Range("D37").Value = "'2.89"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.45%  '

$ws.Range("E38").Value = '  +3.59%  '

$ws.Range("E39").Value = '  +4.26%  '

$ws.Range("D40").Value = "'1.78"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.53%  '

$ws.Range("D41").Value = "'3.94"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.60%  '

$ws.Range("E42").Value = '  -1.45%  '

$ws.Range("D43").Value = '1.966.08'
$ws.Range("E43").Value = '  +2.82%  '

$ws.Range("D44").Value = "'0.0282"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.58%  '

$ws.Range("D45").Value = "'18.54"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.86%  '

$ws.Range("D46").Value = "'2.91"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.45%  '

$ws.Range("D47").Value = "'8.97"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +5.27%  '

$ws.Range("D48").Value = '2.697.11'
$ws.Range("E48").Value = '  +0.91%  '

$ws.Range("D49").Value = "'96.50"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.73%  '

$ws.Range("D50").Value = "'66.49"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.67%  '

$ws.Range("D51").Value = "'0.172"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.40%  '
